$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add Port.new value
$ws.Range("G2").Value = 5

# Row 3: add Port.old and Port.new values
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 6

# Row 7 previously held the "FB8 / Port 5 went bad" record, which is removed.
# The old row 8 record (TV4/S3/PYR/MC1...) moves up to become the new row 7.
$ws.Range("A7").Value = "TV4"
$ws.Range("B7").Value = "S3"
$ws.Range("C7").Value = "PYR"
$ws.Range("D7").Value = "MC1"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = "MC3"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 45401
$ws.Range("I7").Value = "Port went bad"

# Row 8 becomes a brand-new log entry, with no Note.
$ws.Range("A8").Value = "TV1"
$ws.Range("B8").Value = "S5"
$ws.Range("C8").Value = "A14"
$ws.Range("D8").Value = "MC1"
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = "MC1"
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 45492
$ws.Range("I8").ClearContents()

# Update the saved selection to match the authored state.
$ws.Range("G3").Select()
